# Update "want to go" counts (column F) across the three data sheets of the
# workbook to reflect a newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1855
$ws1.Range("F4").Value = 141
$ws1.Range("F6").Value = 3131
$ws1.Range("F10").Value = 638
$ws1.Range("F11").Value = 545
$ws1.Range("F13").Value = 391
$ws1.Range("F18").Value = 1625
$ws1.Range("F19").Value = 14
$ws1.Range("F22").Value = 13
$ws1.Range("F28").Value = 6
$ws1.Range("F30").Value = 33
$ws1.Range("F31").Value = 86
$ws1.Range("F32").Value = 3902
$ws1.Range("F36").Value = 1102
$ws1.Range("F38").Value = 1864

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 25

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1855
$ws4.Range("F4").Value = 141
$ws4.Range("F6").Value = 3131
$ws4.Range("F10").Value = 638
$ws4.Range("F11").Value = 545
$ws4.Range("F13").Value = 25
$ws4.Range("F14").Value = 391
$ws4.Range("F19").Value = 1625
$ws4.Range("F20").Value = 14
$ws4.Range("F23").Value = 13
$ws4.Range("F29").Value = 6
$ws4.Range("F31").Value = 33
$ws4.Range("F32").Value = 86
$ws4.Range("F33").Value = 3902
$ws4.Range("F38").Value = 1102
$ws4.Range("F40").Value = 1864
